# Insert a new "testDeleteUser" test-result row into the UserData sheet's
# results table, pushing the previously-empty template row (row 3) down to
# row 4 (keeping its border styling) and filling row 3 with the new data,
# formatted like the existing "testPostUser" row (row 2 / row 1 header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# 1) Push the current (blank) row 3 formatting down into row 4, so row 4
#    keeps the look that used to belong to row 3, before we overwrite row 3.
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# 2) Give row 3 (A:D) the same formatting as the header/data rows above it.
$ws.Range("A1:D1").Copy()
$ws.Range("A3:D3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# 3) Column E of row 3 should match column E of row 2 (data-row style).
$ws.Range("E2").Copy()
$ws.Range("E3").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# 4) Fill in the new test data for the "testDeleteUser" test case.
$ws.Range("A3").Value = "testDeleteUser"
$ws.Range("B3").Value = "Ayush"
$ws.Range("C3").Value = "Leader"
$ws.Range("D3").Value = "204"
